# "Generate Report for Archive"
#
# Status text update: every cell still showing the old localization status
# ("Ready for handoff") moves to "In Translation" - this affects the
# Overview summary sheet (columns "zh-cn" / "de-de", E2 & F2) as well as
# the per-locale detail sheets' "Status" column (C2) on the "zh-cn" and
# "de-de" tabs.
#
# Alongside that, the "Status" column is narrowed on all three sheets
# (Overview's E & F columns, and column C on the two locale sheets) to fit
# the new, shorter status text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: zh-cn (E2) / de-de (F2) status cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# --- zh-cn sheet: Status column (C2) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus

# --- de-de sheet: Status column (C2) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus

# --- Narrow the status columns to match the shorter text ---
# ColumnWidth is expressed in "characters"; 12.5 is the closest settable
# value to the target stored width on this engine's pixel-quantized grid.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
